# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.896.92"
$ws.Range("E2").Value = "'  +2.78%  "

$ws.Range("D3").Value = "'3.439.79"
$ws.Range("E3").Value = "'  +1.82%  "

$ws.Range("E4").Value = "'  -0.04%  "

$ws.Range("D5").Value = "'580.68"
$ws.Range("E5").Value = "'  +4.58%  "

$ws.Range("D6").Value = "'188.80"
$ws.Range("E6").Value = "'  +8.44%  "

$ws.Range("E7").Value = "'  -0.39%  "

$ws.Range("D8").Value = "'3.431.86"
$ws.Range("E8").Value = "'  +1.81%  "

$ws.Range("E10").Value = "'  -0.98%  "

$ws.Range("E11").Value = "'  +1.39%  "

$ws.Range("D12").Value = "'57.21"
$ws.Range("E12").Value = "'  +6.86%  "

$ws.Range("E13").Value = "'  -0.70%  "

$ws.Range("D14").Value = "'9.45"
$ws.Range("E14").Value = "'  +3.38%  "

$ws.Range("D15").Value = "'3.986.47"
$ws.Range("E15").Value = "'  +1.78%  "

$ws.Range("D16").Value = "'18.82"
$ws.Range("E16").Value = "'  +2.78%  "

$ws.Range("D17").Value = "'3.439.58"
$ws.Range("E17").Value = "'  +1.45%  "

$ws.Range("D18").Value = "'66.836.12"
$ws.Range("E18").Value = "'  +2.85%  "

$ws.Range("E19").Value = "'  +0.32%  "

$ws.Range("E20").Value = "'  +2.16%  "

$ws.Range("E21").Value = "'  +2.83%  "

$ws.Range("D22").Value = "'476.88"
$ws.Range("E22").Value = "'  +4.30%  "

$ws.Range("D23").Value = "'5.41"
$ws.Range("E23").Value = "'  +11.14%  "

$ws.Range("D24").Value = "'17.10"
$ws.Range("E24").Value = "'  +20.79%  "

$ws.Range("E25").Value = "'  +6.90%  "

$ws.Range("D26").Value = "'89.18"
$ws.Range("E26").Value = "'  +1.69%  "

$ws.Range("E27").Value = "'  +3.44%  "

$ws.Range("D28").Value = "'10.93"
$ws.Range("E28").Value = "'  +2.47%  "

$ws.Range("D29").Value = "'8.99"
$ws.Range("E29").Value = "'  +3.54%  "

$ws.Range("D30").Value = "'31.08"
$ws.Range("E30").Value = "'  +0.08%  "

$ws.Range("E31").Value = "'  +13.81%  "

$ws.Range("B32").Value = "'OKB"
$ws.Range("C32").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "'64.70"
$ws.Range("E32").Value = "'  +2.39%  "

$ws.Range("B33").Value = "'Bittensor"
$ws.Range("C33").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'598.96"
$ws.Range("E33").Value = "'  +3.89%  "

$ws.Range("D34").Value = "'11.76"
$ws.Range("E34").Value = "'  +2.75%  "

$ws.Range("E35").Value = "'  +3.83%  "

$ws.Range("E36").Value = "'  -0.04%  "

$ws.Range("E37").Value = "'  +3.37%  "

$ws.Range("D38").Value = "'37.05"
$ws.Range("E38").Value = "'  +4.07%  "

$ws.Range("E39").Value = "'  +4.71%  "

$ws.Range("E40").Value = "'  -4.17%  "

$ws.Range("D41").Value = "'0.0₃0751"
$ws.Range("E41").Value = "'  +1.85%  "

$ws.Range("D42").Value = "'3.201.36"
$ws.Range("E42").Value = "'  +3.21%  "

$ws.Range("D43").Value = "'2.90"
$ws.Range("E43").Value = "'  +5.65%  "

$ws.Range("E44").Value = "'  +3.45%  "

$ws.Range("E45").Value = "'  +5.78%  "

$ws.Range("E46").Value = "'  +1.52%  "

$ws.Range("D47").Value = "'2.72"
$ws.Range("E47").Value = "'  +20.96%  "

$ws.Range("E48").Value = "'  +0.68%  "

$ws.Range("E49").Value = "'  +0.03%  "

$ws.Range("B50").Value = "'LidoDAOToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "'3.18"
$ws.Range("E50").Value = "'  +5.32%  "

$ws.Range("B51").Value = "'THORChain"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.57"
$ws.Range("E51").Value = "'  +3.13%  "
